$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 12, shifting existing
# rows 12..43 down to 13..44 (and extending the used range to row 44).
$ws.Rows("12:12").Insert()

# Populate the newly inserted row 12 with the new weekly record.
$ws.Range("A12").Value = 4
$ws.Range("B12").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C12").Value = "Los Lagos"
$ws.Range("D12").Value = 44624
$ws.Range("E12").Value = 10
$ws.Range("F12").Value = 100112043
$ws.Range("G12").Value = "Pepino dulce"
$ws.Range("H12").Value = "Cultivar IV Región"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 150
$ws.Range("K12").Value = 18000
$ws.Range("L12").Value = 18000
$ws.Range("M12").Value = 18000
$ws.Range("N12").Value = "$/bandeja 18 kilos"
$ws.Range("O12").Value = "Provincia de Limarí"
$ws.Range("P12").Value = 1000
$ws.Range("Q12").Value = 18
$ws.Range("R12").Value = "Hortaliza"
